$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.369.80'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.188.45'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.22'
$ws.Range("E5").Value = '  +1.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.94'
$ws.Range("E6").Value = '  +3.05%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.189.93'
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +2.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.87'
$ws.Range("E11").Value = '  -4.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.510'
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("E13").Value = '  -2.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.94'
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.711.50'
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.460.27'
$ws.Range("E16").Value = '  +0.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.40'
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.188.08'
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '513.37'
$ws.Range("E20").Value = '  +0.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.49'
$ws.Range("E21").Value = '  -3.19%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("E23").Value = '  +2.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.90'
$ws.Range("E24").Value = '  -2.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.69'
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.00'
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.20'
$ws.Range("E28").Value = '  -0.36%  '
$ws.Range("E29").Value = '  +7.09%  '
$ws.Range("E30").Value = '  +7.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.04'
$ws.Range("E31").Value = '  +5.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.11'
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("E33").Value = '  -1.37%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("E35").Value = '  -0.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '508.84'
$ws.Range("E36").Value = '  +5.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.76'
$ws.Range("E37").Value = '  -1.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0897'
$ws.Range("E38").Value = '  -2.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0423'
$ws.Range("E39").Value = '  -0.24%  '
$ws.Range("E40").Value = '  +6.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.87'
$ws.Range("E41").Value = '  -1.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.301'
$ws.Range("E42").Value = '  +4.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.85'
$ws.Range("E43").Value = '  -5.08%  '
$ws.Range("E44").Value = '  +5.17%  '
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.852.98'
$ws.Range("E46").Value = '  -5.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.43'
$ws.Range("E47").Value = '  -2.29%  '
$ws.Range("E48").Value = '  +5.26%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("E50").Value = '  +0.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.58'
$ws.Range("E51").Value = '  +4.84%  '
